# Fruta / hortaliza, semanal
# Weekly update: insert a new observation row for Mango (Vega Central Mapocho de
# Santiago) at row 630, pushing the existing rows 630..758 down to 631..759.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 630 (shifts rows 630:758 down to 631:759)
$ws.Rows.Item(630).Insert()

# Populate the newly inserted row with the new weekly price observation
$ws.Range("A630").Value = 9
$ws.Range("B630").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C630").Value = "Metropolitana"
$ws.Range("D630").Value = 45209
$ws.Range("E630").Value = 13
$ws.Range("F630").Value = "Fruta"
$ws.Range("G630").Value = 100108
$ws.Range("H630").Value = "Tropicales y subtropicales"
$ws.Range("I630").Value = 100108002
$ws.Range("J630").Value = "Mango"
$ws.Range("K630").Value = "Sin especificar"
$ws.Range("L630").Value = "Primera"
$ws.Range("M630").Value = 410
$ws.Range("N630").Value = 10000
$ws.Range("O630").Value = 10000
$ws.Range("P630").Value = 10000
$ws.Range("Q630").Value = "`$/bandeja 4 kilos"
$ws.Range("R630").Value = "Brasil"
$ws.Range("S630").Value = 2500
$ws.Range("T630").Value = 4
